$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New StatQuery text (replaces the old StatQuery text previously shared by C2, C3, C4)
$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Giant Schnauzer']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Column D no longer needs to be as wide as before now that its neighboring
# column content changed; narrow it down (re-fit) to match the new layout.
$ws.Columns.Item(4).ColumnWidth = 44.25

# Move/restore the active selection to B4 (author's final cursor position)
$ws.Range("B4").Select()
